$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E1 text value from "c" to "cc"
$ws.Range("E1").Value = "cc"

# Update C4 value from 0 to 5
$ws.Range("C4").Value = 5

# Update selection to E12
$ws.Range("E12").Select()
